$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Apply header style (bold, bordered, centered) to new row-index cells A20:A23,
# matching the style already used for A2:A19.
$ws.Range("A16:A19").Copy() | Out-Null
$ws.Range("A20:A23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(1,2).Value = 0
$ws.Cells.Item(1,3).Value = 1
$ws.Cells.Item(1,4).Value = 2
$ws.Cells.Item(1,5).Value = 3
$ws.Cells.Item(1,6).Value = 4
$ws.Cells.Item(1,7).Value = 5
$ws.Cells.Item(1,8).Value = 6
$ws.Cells.Item(1,9).Value = 7
$ws.Cells.Item(1,10).Value = 8
$ws.Cells.Item(1,11).Value = 9
$ws.Cells.Item(1,12).Value = 10
$ws.Cells.Item(1,13).Value = 11
$ws.Cells.Item(1,14).Value = 12
$ws.Cells.Item(1,15).Value = 13
$ws.Cells.Item(1,16).Value = 14
$ws.Cells.Item(1,17).Value = 15
$ws.Cells.Item(1,18).Value = 16
$ws.Cells.Item(1,19).Value = 17
$ws.Cells.Item(1,20).Value = 18
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "HKL"
$ws.Cells.Item(2,3).Value = "[2, 0, 0]"
$ws.Cells.Item(2,4).Value = "[2, 2, 0]"
$ws.Cells.Item(2,5).Value = "[4, 0, 0]"
$ws.Cells.Item(2,6).Value = "[2, 1, 1]"
$ws.Cells.Item(2,7).Value = "[3, 2, 1]"
$ws.Cells.Item(2,8).Value = "[3, 1, 0]"
$ws.Cells.Item(2,9).Value = "[2, 2, 2]"
$ws.Cells.Item(2,10).Value = "[1, 1, 0]"
$ws.Cells.Item(2,11).Value = "1Pair-A"
$ws.Cells.Item(2,12).Value = "1Pair-B"
$ws.Cells.Item(2,13).Value = "2Pairs-A"
$ws.Cells.Item(2,14).Value = "2Pairs-B"
$ws.Cells.Item(2,15).Value = "3Pairs-A"
$ws.Cells.Item(2,16).Value = "3Pairs-B"
$ws.Cells.Item(2,17).Value = "3Pairs-C"
$ws.Cells.Item(2,18).Value = "4Pairs"
$ws.Cells.Item(2,19).Value = "5A4F"
$ws.Cells.Item(2,20).Value = "MaxUnique"
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "BT8Hex_2.5"
$ws.Cells.Item(3,3).Value = 0.9603982863342367
$ws.Cells.Item(3,4).Value = 1.033295320830496
$ws.Cells.Item(3,5).Value = 0.9603982863342367
$ws.Cells.Item(3,6).Value = 0.9568466131739304
$ws.Cells.Item(3,7).Value = 0.9952719105738961
$ws.Cells.Item(3,8).Value = 1.001687357152962
$ws.Cells.Item(3,9).Value = 0.9480679269782257
$ws.Cells.Item(3,10).Value = 1.033295320830496
$ws.Cells.Item(3,11).Value = 1.033295320830496
$ws.Cells.Item(3,12).Value = 0.9568466131739304
$ws.Cells.Item(3,13).Value = 0.9586224497540836
$ws.Cells.Item(3,14).Value = 0.9586224497540836
$ws.Cells.Item(3,15).Value = 0.9729774188870431
$ws.Cells.Item(3,16).Value = 0.9835134067795545
$ws.Cells.Item(3,17).Value = 0.9835134067795545
$ws.Cells.Item(3,18).Value = 0.9959588852922899
$ws.Cells.Item(3,19).Value = 0.9959588852922899
$ws.Cells.Item(3,20).Value = 0.9825945691739579
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "BT8Hex_5"
$ws.Cells.Item(4,3).Value = 0.9490389450926037
$ws.Cells.Item(4,4).Value = 1.078802738517202
$ws.Cells.Item(4,5).Value = 0.9490389450926037
$ws.Cells.Item(4,6).Value = 0.9374798231982281
$ws.Cells.Item(4,7).Value = 1.003098938229908
$ws.Cells.Item(4,8).Value = 0.972347547405844
$ws.Cells.Item(4,9).Value = 0.9194045203581911
$ws.Cells.Item(4,10).Value = 1.078802738517202
$ws.Cells.Item(4,11).Value = 1.078802738517202
$ws.Cells.Item(4,12).Value = 0.9374798231982281
$ws.Cells.Item(4,13).Value = 0.9432593841454159
$ws.Cells.Item(4,14).Value = 0.9432593841454159
$ws.Cells.Item(4,15).Value = 0.9529554385655586
$ws.Cells.Item(4,16).Value = 0.9884405022693444
$ws.Cells.Item(4,17).Value = 0.9884405022693444
$ws.Cells.Item(4,18).Value = 1.011031061331309
$ws.Cells.Item(4,19).Value = 1.011031061331309
$ws.Cells.Item(4,20).Value = 0.9766954188003293
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "BT8Hex_10"
$ws.Cells.Item(5,3).Value = 0.5340260645028635
$ws.Cells.Item(5,4).Value = 1.245390163545357
$ws.Cells.Item(5,5).Value = 0.5340260645028635
$ws.Cells.Item(5,6).Value = 0.8115824054734082
$ws.Cells.Item(5,7).Value = 1.112894864222037
$ws.Cells.Item(5,8).Value = 0.887117413513217
$ws.Cells.Item(5,9).Value = 1.001665038562908
$ws.Cells.Item(5,10).Value = 1.245390163545357
$ws.Cells.Item(5,11).Value = 1.245390163545357
$ws.Cells.Item(5,12).Value = 0.8115824054734082
$ws.Cells.Item(5,13).Value = 0.6728042349881358
$ws.Cells.Item(5,14).Value = 0.6728042349881358
$ws.Cells.Item(5,15).Value = 0.7442419611631629
$ws.Cells.Item(5,16).Value = 0.8636662111738763
$ws.Cells.Item(5,17).Value = 0.8636662111738763
$ws.Cells.Item(5,18).Value = 0.9590971992667465
$ws.Cells.Item(5,19).Value = 0.9590971992667465
$ws.Cells.Item(5,20).Value = 0.9321126583032985
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "BT8Hex_15"
$ws.Cells.Item(6,3).Value = 0.2524654590984607
$ws.Cells.Item(6,4).Value = 1.53919425166635
$ws.Cells.Item(6,5).Value = 0.2524654590984607
$ws.Cells.Item(6,6).Value = 0.5076075292482196
$ws.Cells.Item(6,7).Value = 1.186726939150263
$ws.Cells.Item(6,8).Value = 0.908075979937792
$ws.Cells.Item(6,9).Value = 0.7403039790582365
$ws.Cells.Item(6,10).Value = 1.53919425166635
$ws.Cells.Item(6,11).Value = 1.53919425166635
$ws.Cells.Item(6,12).Value = 0.5076075292482196
$ws.Cells.Item(6,13).Value = 0.3800364941733402
$ws.Cells.Item(6,14).Value = 0.3800364941733402
$ws.Cells.Item(6,15).Value = 0.5560496560948242
$ws.Cells.Item(6,16).Value = 0.7664224133376768
$ws.Cells.Item(6,17).Value = 0.7664224133376768
$ws.Cells.Item(6,18).Value = 0.9596153729198451
$ws.Cells.Item(6,19).Value = 0.9596153729198451
$ws.Cells.Item(6,20).Value = 0.8557290230265536
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "Spiral2.5"
$ws.Cells.Item(7,3).Value = 0.999208457531287
$ws.Cells.Item(7,4).Value = 1.003455079705105
$ws.Cells.Item(7,5).Value = 0.999208457531287
$ws.Cells.Item(7,6).Value = 0.9988543346555983
$ws.Cells.Item(7,7).Value = 0.9985097368416292
$ws.Cells.Item(7,8).Value = 0.9990965098792857
$ws.Cells.Item(7,9).Value = 0.998797493702943
$ws.Cells.Item(7,10).Value = 1.003455079705105
$ws.Cells.Item(7,11).Value = 1.003455079705105
$ws.Cells.Item(7,12).Value = 0.9988543346555983
$ws.Cells.Item(7,13).Value = 0.9990313960934426
$ws.Cells.Item(7,14).Value = 0.9990313960934426
$ws.Cells.Item(7,15).Value = 0.9990531006887237
$ws.Cells.Item(7,16).Value = 1.00050595729733
$ws.Cells.Item(7,17).Value = 1.00050595729733
$ws.Cells.Item(7,18).Value = 1.001243237899274
$ws.Cells.Item(7,19).Value = 1.001243237899274
$ws.Cells.Item(7,20).Value = 0.9996536020526414
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "Spiral5"
$ws.Cells.Item(8,3).Value = 0.9983368492440228
$ws.Cells.Item(8,4).Value = 1.016783358394947
$ws.Cells.Item(8,5).Value = 0.9983368492440228
$ws.Cells.Item(8,6).Value = 0.9969831909762562
$ws.Cells.Item(8,7).Value = 0.9963476550026472
$ws.Cells.Item(8,8).Value = 0.9996274053393231
$ws.Cells.Item(8,9).Value = 1.005829307328
$ws.Cells.Item(8,10).Value = 1.016783358394947
$ws.Cells.Item(8,11).Value = 1.016783358394947
$ws.Cells.Item(8,12).Value = 0.9969831909762562
$ws.Cells.Item(8,13).Value = 0.9976600201101395
$ws.Cells.Item(8,14).Value = 0.9976600201101395
$ws.Cells.Item(8,15).Value = 0.9983158151865341
$ws.Cells.Item(8,16).Value = 1.004034466205075
$ws.Cells.Item(8,17).Value = 1.004034466205075
$ws.Cells.Item(8,18).Value = 1.007221689252543
$ws.Cells.Item(8,19).Value = 1.007221689252543
$ws.Cells.Item(8,20).Value = 1.002317961047533
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "Spiral7.5"
$ws.Cells.Item(9,3).Value = 1.002521294339786
$ws.Cells.Item(9,4).Value = 1.047408718375914
$ws.Cells.Item(9,5).Value = 1.002521294339786
$ws.Cells.Item(9,6).Value = 0.9997040670494584
$ws.Cells.Item(9,7).Value = 0.9933300795949923
$ws.Cells.Item(9,8).Value = 0.9994280159427856
$ws.Cells.Item(9,9).Value = 0.9903562844195494
$ws.Cells.Item(9,10).Value = 1.047408718375914
$ws.Cells.Item(9,11).Value = 1.047408718375914
$ws.Cells.Item(9,12).Value = 0.9997040670494584
$ws.Cells.Item(9,13).Value = 1.001112680694622
$ws.Cells.Item(9,14).Value = 1.001112680694622
$ws.Cells.Item(9,15).Value = 1.000551125777343
$ws.Cells.Item(9,16).Value = 1.016544693255053
$ws.Cells.Item(9,17).Value = 1.016544693255053
$ws.Cells.Item(9,18).Value = 1.024260699535268
$ws.Cells.Item(9,19).Value = 1.024260699535268
$ws.Cells.Item(9,20).Value = 1.005458076620414
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "Spiral10"
$ws.Cells.Item(10,3).Value = 0.9920446302045643
$ws.Cells.Item(10,4).Value = 1.071051514714167
$ws.Cells.Item(10,5).Value = 0.9920446302045643
$ws.Cells.Item(10,6).Value = 0.9950890750067719
$ws.Cells.Item(10,7).Value = 0.9954105318570978
$ws.Cells.Item(10,8).Value = 0.9917994037693324
$ws.Cells.Item(10,9).Value = 0.9929555741231759
$ws.Cells.Item(10,10).Value = 1.071051514714167
$ws.Cells.Item(10,11).Value = 1.071051514714167
$ws.Cells.Item(10,12).Value = 0.9950890750067719
$ws.Cells.Item(10,13).Value = 0.9935668526056681
$ws.Cells.Item(10,14).Value = 0.9935668526056681
$ws.Cells.Item(10,15).Value = 0.9929777029935561
$ws.Cells.Item(10,16).Value = 1.019395073308501
$ws.Cells.Item(10,17).Value = 1.019395073308501
$ws.Cells.Item(10,18).Value = 1.032309183659917
$ws.Cells.Item(10,19).Value = 1.032309183659917
$ws.Cells.Item(10,20).Value = 1.006391788279185
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "Spiral15"
$ws.Cells.Item(11,3).Value = 1.002621428100762
$ws.Cells.Item(11,4).Value = 1.250814630877674
$ws.Cells.Item(11,5).Value = 1.002621428100762
$ws.Cells.Item(11,6).Value = 0.9843042375330987
$ws.Cells.Item(11,7).Value = 0.9894431392254502
$ws.Cells.Item(11,8).Value = 0.9789800100535136
$ws.Cells.Item(11,9).Value = 0.8355945338166691
$ws.Cells.Item(11,10).Value = 1.250814630877674
$ws.Cells.Item(11,11).Value = 1.250814630877674
$ws.Cells.Item(11,12).Value = 0.9843042375330987
$ws.Cells.Item(11,13).Value = 0.9934628328169306
$ws.Cells.Item(11,14).Value = 0.9934628328169306
$ws.Cells.Item(11,15).Value = 0.988635225229125
$ws.Cells.Item(11,16).Value = 1.079246765503845
$ws.Cells.Item(11,17).Value = 1.079246765503845
$ws.Cells.Item(11,18).Value = 1.122138731847303
$ws.Cells.Item(11,19).Value = 1.122138731847303
$ws.Cells.Item(11,20).Value = 1.006959663267861
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "OffsetF45"
$ws.Cells.Item(12,3).Value = 0.02816246208286063
$ws.Cells.Item(12,4).Value = 1.947503822981913
$ws.Cells.Item(12,5).Value = 0.02816246208286063
$ws.Cells.Item(12,6).Value = 0.893658514633978
$ws.Cells.Item(12,7).Value = 1.444407193392876
$ws.Cells.Item(12,8).Value = 1.088287096789949
$ws.Cells.Item(12,9).Value = 0.008384910052905486
$ws.Cells.Item(12,10).Value = 1.947503822981913
$ws.Cells.Item(12,11).Value = 1.947503822981913
$ws.Cells.Item(12,12).Value = 0.893658514633978
$ws.Cells.Item(12,13).Value = 0.4609104883584194
$ws.Cells.Item(12,14).Value = 0.4609104883584194
$ws.Cells.Item(12,15).Value = 0.6700360245022626
$ws.Cells.Item(12,16).Value = 0.956441599899584
$ws.Cells.Item(12,17).Value = 0.956441599899584
$ws.Cells.Item(12,18).Value = 1.204207155670166
$ws.Cells.Item(12,19).Value = 1.204207155670166
$ws.Cells.Item(12,20).Value = 0.9017339999890804
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "OffsetA45"
$ws.Cells.Item(13,3).Value = 1.978384880534225
$ws.Cells.Item(13,4).Value = 0.8299688759072553
$ws.Cells.Item(13,5).Value = 1.978384880534225
$ws.Cells.Item(13,6).Value = 1.050809254152312
$ws.Cells.Item(13,7).Value = 1.105508623078699
$ws.Cells.Item(13,8).Value = 1.324511104887578
$ws.Cells.Item(13,9).Value = 1.97000817876325
$ws.Cells.Item(13,10).Value = 0.8299688759072553
$ws.Cells.Item(13,11).Value = 0.8299688759072553
$ws.Cells.Item(13,12).Value = 1.050809254152312
$ws.Cells.Item(13,13).Value = 1.514597067343268
$ws.Cells.Item(13,14).Value = 1.514597067343268
$ws.Cells.Item(13,15).Value = 1.451235079858038
$ws.Cells.Item(13,16).Value = 1.286387670197931
$ws.Cells.Item(13,17).Value = 1.286387670197931
$ws.Cells.Item(13,18).Value = 1.172282971625262
$ws.Cells.Item(13,19).Value = 1.172282971625262
$ws.Cells.Item(13,20).Value = 1.376531819553886
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "OffsetFTD"
$ws.Cells.Item(14,3).Value = 1.139695438201505
$ws.Cells.Item(14,4).Value = -0.000422063954842297
$ws.Cells.Item(14,5).Value = 1.139695438201505
$ws.Cells.Item(14,6).Value = 1.453649013070499
$ws.Cells.Item(14,7).Value = 0.5717459883016168
$ws.Cells.Item(14,8).Value = 1.70883991090389
$ws.Cells.Item(14,9).Value = 1.141189664377114
$ws.Cells.Item(14,10).Value = -0.000422063954842297
$ws.Cells.Item(14,11).Value = -0.000422063954842297
$ws.Cells.Item(14,12).Value = 1.453649013070499
$ws.Cells.Item(14,13).Value = 1.296672225636002
$ws.Cells.Item(14,14).Value = 1.296672225636002
$ws.Cells.Item(14,15).Value = 1.434061454058631
$ws.Cells.Item(14,16).Value = 0.864307462439054
$ws.Cells.Item(14,17).Value = 0.864307462439054
$ws.Cells.Item(14,18).Value = 0.64812508084058
$ws.Cells.Item(14,19).Value = 0.64812508084058
$ws.Cells.Item(14,20).Value = 1.002449658483297
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "OffsetATD"
$ws.Cells.Item(15,3).Value = 0.7216725738718084
$ws.Cells.Item(15,4).Value = 0.8100617453708652
$ws.Cells.Item(15,5).Value = 0.7216725738718084
$ws.Cells.Item(15,6).Value = 0.8021127364917493
$ws.Cells.Item(15,7).Value = 1.304452449442555
$ws.Cells.Item(15,8).Value = 0.1598485097874724
$ws.Cells.Item(15,9).Value = 0.1338779895585501
$ws.Cells.Item(15,10).Value = 0.8100617453708652
$ws.Cells.Item(15,11).Value = 0.8100617453708652
$ws.Cells.Item(15,12).Value = 0.8021127364917493
$ws.Cells.Item(15,13).Value = 0.7618926551817788
$ws.Cells.Item(15,14).Value = 0.7618926551817788
$ws.Cells.Item(15,15).Value = 0.5612112733836767
$ws.Cells.Item(15,16).Value = 0.7779490185781409
$ws.Cells.Item(15,17).Value = 0.7779490185781409
$ws.Cells.Item(15,18).Value = 0.785977200276322
$ws.Cells.Item(15,19).Value = 0.785977200276322
$ws.Cells.Item(15,20).Value = 0.6553376674205
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "Holden2.5"
$ws.Cells.Item(16,3).Value = 0.4896915439846395
$ws.Cells.Item(16,4).Value = 1.563726112658459
$ws.Cells.Item(16,5).Value = 0.4896915439846395
$ws.Cells.Item(16,6).Value = 0.3600815523740016
$ws.Cells.Item(16,7).Value = 0.9049710619826785
$ws.Cells.Item(16,8).Value = 0.98851176542111
$ws.Cells.Item(16,9).Value = 0.2006914790345002
$ws.Cells.Item(16,10).Value = 1.563726112658459
$ws.Cells.Item(16,11).Value = 1.563726112658459
$ws.Cells.Item(16,12).Value = 0.3600815523740016
$ws.Cells.Item(16,13).Value = 0.4248865481793206
$ws.Cells.Item(16,14).Value = 0.4248865481793206
$ws.Cells.Item(16,15).Value = 0.6127616205932503
$ws.Cells.Item(16,16).Value = 0.8044997363390335
$ws.Cells.Item(16,17).Value = 0.8044997363390335
$ws.Cells.Item(16,18).Value = 0.9943063304188899
$ws.Cells.Item(16,19).Value = 0.9943063304188899
$ws.Cells.Item(16,20).Value = 0.7512789192425648
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "Holden5"
$ws.Cells.Item(17,3).Value = 0.7401734415562236
$ws.Cells.Item(17,4).Value = 1.332717961315786
$ws.Cells.Item(17,5).Value = 0.7401734415562236
$ws.Cells.Item(17,6).Value = 0.5322704974433264
$ws.Cells.Item(17,7).Value = 0.8911587992759364
$ws.Cells.Item(17,8).Value = 1.041398673533703
$ws.Cells.Item(17,9).Value = 0.4014102716328738
$ws.Cells.Item(17,10).Value = 1.332717961315786
$ws.Cells.Item(17,11).Value = 1.332717961315786
$ws.Cells.Item(17,12).Value = 0.5322704974433264
$ws.Cells.Item(17,13).Value = 0.636221969499775
$ws.Cells.Item(17,14).Value = 0.636221969499775
$ws.Cells.Item(17,15).Value = 0.7712808708444175
$ws.Cells.Item(17,16).Value = 0.868387300105112
$ws.Cells.Item(17,17).Value = 0.8683873001051118
$ws.Cells.Item(17,18).Value = 0.9844699654077802
$ws.Cells.Item(17,19).Value = 0.9844699654077802
$ws.Cells.Item(17,20).Value = 0.823188274126308
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = "Holden10"
$ws.Cells.Item(18,3).Value = 1.223818666748404
$ws.Cells.Item(18,4).Value = 0.8998328563097706
$ws.Cells.Item(18,5).Value = 1.223818666748404
$ws.Cells.Item(18,6).Value = 0.8734191763370657
$ws.Cells.Item(18,7).Value = 0.8599796921775044
$ws.Cells.Item(18,8).Value = 1.144116811400104
$ws.Cells.Item(18,9).Value = 0.8122374857920637
$ws.Cells.Item(18,10).Value = 0.8998328563097706
$ws.Cells.Item(18,11).Value = 0.8998328563097706
$ws.Cells.Item(18,12).Value = 0.8734191763370657
$ws.Cells.Item(18,13).Value = 1.048618921542735
$ws.Cells.Item(18,14).Value = 1.048618921542735
$ws.Cells.Item(18,15).Value = 1.080451551495191
$ws.Cells.Item(18,16).Value = 0.9990235664650801
$ws.Cells.Item(18,17).Value = 0.9990235664650801
$ws.Cells.Item(18,18).Value = 0.9742258889262527
$ws.Cells.Item(18,19).Value = 0.9742258889262527
$ws.Cells.Item(18,20).Value = 0.9689007814608188
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = "Holden15"
$ws.Cells.Item(19,3).Value = 1.208742345840154
$ws.Cells.Item(19,4).Value = 0.7923945562693139
$ws.Cells.Item(19,5).Value = 1.208742345840154
$ws.Cells.Item(19,6).Value = 0.846781111899922
$ws.Cells.Item(19,7).Value = 0.8734813319885301
$ws.Cells.Item(19,8).Value = 1.118880797814664
$ws.Cells.Item(19,9).Value = 0.8520880594203418
$ws.Cells.Item(19,10).Value = 0.7923945562693139
$ws.Cells.Item(19,11).Value = 0.7923945562693139
$ws.Cells.Item(19,12).Value = 0.846781111899922
$ws.Cells.Item(19,13).Value = 1.027761728870038
$ws.Cells.Item(19,14).Value = 1.027761728870038
$ws.Cells.Item(19,15).Value = 1.05813475185158
$ws.Cells.Item(19,16).Value = 0.9493060046697966
$ws.Cells.Item(19,17).Value = 0.9493060046697966
$ws.Cells.Item(19,18).Value = 0.910078142569676
$ws.Cells.Item(19,19).Value = 0.910078142569676
$ws.Cells.Item(19,20).Value = 0.9487280338721543
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = "HexGrid-90degTilt2.5degRes"
$ws.Cells.Item(20,3).Value = 0.9998098638866774
$ws.Cells.Item(20,4).Value = 0.9915557112424485
$ws.Cells.Item(20,5).Value = 0.9998098638866774
$ws.Cells.Item(20,6).Value = 0.9997025313974955
$ws.Cells.Item(20,7).Value = 1.000517729476282
$ws.Cells.Item(20,8).Value = 1.001822993843897
$ws.Cells.Item(20,9).Value = 1.000972930761219
$ws.Cells.Item(20,10).Value = 0.9915557112424485
$ws.Cells.Item(20,11).Value = 0.9915557112424485
$ws.Cells.Item(20,12).Value = 0.9997025313974955
$ws.Cells.Item(20,13).Value = 0.9997561976420865
$ws.Cells.Item(20,14).Value = 0.9997561976420865
$ws.Cells.Item(20,15).Value = 1.000445129709357
$ws.Cells.Item(20,16).Value = 0.9970227021755406
$ws.Cells.Item(20,17).Value = 0.9970227021755403
$ws.Cells.Item(20,18).Value = 0.9956559544422674
$ws.Cells.Item(20,19).Value = 0.9956559544422674
$ws.Cells.Item(20,20).Value = 0.9990636267680033
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = "HexGrid-90degTilt5degRes"
$ws.Cells.Item(21,3).Value = 0.9965996471491063
$ws.Cells.Item(21,4).Value = 1.032878453532439
$ws.Cells.Item(21,5).Value = 0.9965996471491063
$ws.Cells.Item(21,6).Value = 0.9977199701763002
$ws.Cells.Item(21,7).Value = 1.000169252366526
$ws.Cells.Item(21,8).Value = 0.9989794038845281
$ws.Cells.Item(21,9).Value = 1.002300381124696
$ws.Cells.Item(21,10).Value = 1.032878453532439
$ws.Cells.Item(21,11).Value = 1.032878453532439
$ws.Cells.Item(21,12).Value = 0.9977199701763002
$ws.Cells.Item(21,13).Value = 0.9971598086627033
$ws.Cells.Item(21,14).Value = 0.9971598086627033
$ws.Cells.Item(21,15).Value = 0.9977663404033116
$ws.Cells.Item(21,16).Value = 1.009066023619282
$ws.Cells.Item(21,17).Value = 1.009066023619282
$ws.Cells.Item(21,18).Value = 1.015019131097571
$ws.Cells.Item(21,19).Value = 1.015019131097571
$ws.Cells.Item(21,20).Value = 1.004774518038933
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = "HexGrid-90degTilt10degRes"
$ws.Cells.Item(22,3).Value = 0.9979529816929837
$ws.Cells.Item(22,4).Value = 1.142027234367606
$ws.Cells.Item(22,5).Value = 0.9979529816929837
$ws.Cells.Item(22,6).Value = 0.9548748389948418
$ws.Cells.Item(22,7).Value = 1.033693055258503
$ws.Cells.Item(22,8).Value = 0.9938407630462037
$ws.Cells.Item(22,9).Value = 0.9899218139204938
$ws.Cells.Item(22,10).Value = 1.142027234367606
$ws.Cells.Item(22,11).Value = 1.142027234367606
$ws.Cells.Item(22,12).Value = 0.9548748389948418
$ws.Cells.Item(22,13).Value = 0.9764139103439128
$ws.Cells.Item(22,14).Value = 0.9764139103439128
$ws.Cells.Item(22,15).Value = 0.9822228612446765
$ws.Cells.Item(22,16).Value = 1.031618351685144
$ws.Cells.Item(22,17).Value = 1.031618351685144
$ws.Cells.Item(22,18).Value = 1.05922057235576
$ws.Cells.Item(22,19).Value = 1.05922057235576
$ws.Cells.Item(22,20).Value = 1.018718447880105
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = "HexGrid-90degTilt15degRes"
$ws.Cells.Item(23,3).Value = 0.7111392718243243
$ws.Cells.Item(23,4).Value = 0.8045770884521859
$ws.Cells.Item(23,5).Value = 0.7111392718243243
$ws.Cells.Item(23,6).Value = 1.006288395017029
$ws.Cells.Item(23,7).Value = 1.115851747890229
$ws.Cells.Item(23,8).Value = 1.031067143442572
$ws.Cells.Item(23,9).Value = 1.163636169032363
$ws.Cells.Item(23,10).Value = 0.8045770884521859
$ws.Cells.Item(23,11).Value = 0.8045770884521859
$ws.Cells.Item(23,12).Value = 1.006288395017029
$ws.Cells.Item(23,13).Value = 0.8587138334206766
$ws.Cells.Item(23,14).Value = 0.8587138334206766
$ws.Cells.Item(23,15).Value = 0.9161649367613084
$ws.Cells.Item(23,16).Value = 0.8406682517645131
$ws.Cells.Item(23,17).Value = 0.8406682517645131
$ws.Cells.Item(23,18).Value = 0.8316454609364312
$ws.Cells.Item(23,19).Value = 0.8316454609364312
$ws.Cells.Item(23,20).Value = 0.9720933026097839
